# feat: sheetParser now returns an ageGroup array
# Append the new ageGroup rows (7 -> 60/70/100) to the "Mujeres" sheet
# (the active sheet) just below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 60

$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 70

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 100

# Leave the selection where the user last typed, matching the saved view.
[void]$ws.Range("B4").Select()
